# CE09OSPM added eng bar codes and bogey for D4 recovery date
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Asset_Cal_Info")

# Add the new "eng" sensor barcode for the STC row (row 31)
$ws.Range("E31").Style = "Normal"
$ws.Range("E31").Value = "OL000373"

# The RTE000000 placeholder row (row 32) is no longer needed - remove it entirely,
# shifting the WFPENG000 row (formerly row 33) up to row 32.
$ws.Rows("32").Delete()

# The Asset_Cal_Info data range shrank by one row (400 -> 399); update the
# stale _FilterDatabase_* defined names that still pointed at row 400.
$staleNames = @(
    "_FilterDatabase_0_0_0_0_0_0",
    "_FilterDatabase_0_0_0_0_0_0_0_0",
    "_FilterDatabase_0_0_0_0_1",
    "_FilterDatabase_0_0_1",
    "_FilterDatabase_2"
)
foreach ($nm in $staleNames) {
    $wb.Names.Item($nm).RefersTo = "=Asset_Cal_Info!`$A`$1:`$H`$399"
}
